$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new stop word "yg" into the list, in the next empty row (A52)
$ws.Range("A52").Value = "yg"
